# "removed delivery fees with auto import"
#
# The source sheet has a "Delivery Fees" column (H) that is no longer wanted.
# Select the whole column (so the resulting selection matches what Excel
# leaves behind after a column delete) and remove it, shifting every column
# to its right (Markup Delivery Fees, Remark, Full Address, Schedule Date,
# Delivery Type, Payment Method) one place to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Columns("H:H").Select()
[void]$ws.Columns("H:H").Delete()
